$wb = $excel.ActiveWorkbook

# The mail log lives on the "Logs" worksheet.
$ws = $wb.Worksheets.Item("Logs")

# Append the new incoming mail as row 43.
$ws.Range("A43").Value = "Klacht over levering"
$ws.Range("B43").Value = "mailmind.test@zohomail.eu"
$ws.Range("C43").Value = "Mijn bestelling is incompleet geleverd. Graag hoor ik hoe dit wordt opgelost."
$ws.Range("D43").Value = "Bestelling / Levering"
$ws.Range("F43").Value = "2025-06-24 22:11:51"
$ws.Range("G43").Value = "Nee"

# The "Dashboard" sheet keeps a per-category count used by the chart;
# "Bestelling / Levering" just gained one more entry.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B5").Value = 5

# Extend the conditional-formatting ranges so the new row is covered too.
$catFcs = $ws.Range("D2:D42").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D43"))
}

$answeredFcs = $ws.Range("G2:G42").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($ws.Range("G2:G43"))
}
